# Update with new biosteam results
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LCA")

# --- Shared string / label text updates ---
$ws.Range("A29").Value = "item28 [tonne*km]"
$ws.Range("A48").Value = "E_item [kWh]"

# --- GlobalWarming block (rows 2-24): column F ratios recalculated ---
$ws.Range("F2").Value = 0.03049262772909476
$ws.Range("F3").Value = 0.03049262772909476
$ws.Range("F4").Value = 0.3630074729654137
$ws.Range("F5").Value = 0.3630074729654137
$ws.Range("F6").Value = 0.001054450278613821
$ws.Range("F7").Value = 0.004243430824932252
$ws.Range("F8").Value = 0.01332932602379867
$ws.Range("F9").Value = 0.01862720712734474
$ws.Range("F10").Value = 0.0009304131219783373
$ws.Range("F11").Value = 0.00001730527688581152
$ws.Range("F12").Value = 0.0009477183988641488
$ws.Range("F13").Value = 0.009219237408645428
$ws.Range("F14").Value = 0.0002612117265782871
$ws.Range("F15").Value = 0.009480449135223716
$ws.Range("F16").Value = 0.009534996289891535
$ws.Range("F17").Value = 0.009534996289891535
$ws.Range("F18").Value = 0.01827944297199172
$ws.Range("F19").Value = 0.01827944297199172
$ws.Range("F20").Value = 0.04111023618687183
$ws.Range("F22").Value = 0.04111023618687183
$ws.Range("F23").Value = 0.01797271126070825
$ws.Range("F24").Value = 0.01797271126070825
$ws.Range("E25").Value = 6508130.482

# --- Stream contribution block (rows 35-44) ---
$ws.Range("B35").Value = 405530.7657731167
$ws.Range("C35").Value = 11354861.44164727
$ws.Range("D35").Value = 2.127136878583531

$ws.Range("D36").Value = 0.004862779102746523

$ws.Range("B37").Value = 768085.0074808991
$ws.Range("C37").Value = -2306774.376275816
$ws.Range("D37").Value = -0.4321342776011864

$ws.Range("D38").Value = -0.07613945131042262

$ws.Range("B39").Value = 453682.9082841854
$ws.Range("C39").Value = -2449887.704734601
$ws.Range("D39").Value = -0.4589440841625387

$ws.Range("D40").Value = -0.1086100333612747

$ws.Range("D41").Value = -0.01160934631756385

$ws.Range("C42").Value = -94571.20002396787
$ws.Range("D42").Value = -0.01771627846422213

$ws.Range("B43").Value = 29246.43426704999
$ws.Range("D43").Value = -0.02684618646906903

$ws.Range("C44").Value = 5338096.272022005
